$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jRows = @(14,23,27,35,547,559,562,563,568,570,571,574,577,579,580,582,584,587,588,590,592,594,595,597,598,600,601,603,605,606,607,608,609,610,612,614,616,618,620,621,622,625,630,631,633,635,637,638,639,641,642,644,646,647,651,655,656,658,659,665,669,671,673,675,677,679,681,683,685,686,689,691,693,695,696,697,698,699,700,704,707,709,710,712,713,715,716,718,720,721,725,726,729,730,731,733,736,822,824,825,826,827,828,829,830,831,832,833,834,835,836,837,838,839,840,841,842,843,844,845,846,847,848,857,859,860,862,864,867,871,881,884,886,893,896,941,964,970,972,973,985,998,999,1045,1048,1068,1080,1082,1085,1089)
foreach ($r in $jRows) {
    $ws.Range("J$r").Value = "vide"
}

$fRows = @(562,570,621,646,702,723,730,823)
foreach ($r in $fRows) {
    $ws.Range("F$r").Value = "vide"
}
